# Swap the content of row 4 and row 6 (columns A, B, D, E, F, G, H, S),
# and move the "Publik kommentar" (AC) and "Biotop-beskrivning" (AI)
# comments from row 6 to row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture current (row 4 / row 6) values for the columns that swap ---
# NOTE: use Value2 (not Value) - Value is a parameterized COM property
# that this host does not marshal correctly when read without arguments.
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$d4 = $ws.Range("D4").Value2
$e4 = $ws.Range("E4").Value2
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$h4 = $ws.Range("H4").Value2
$s4 = $ws.Range("S4").Value2

$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2
$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2
$f6 = $ws.Range("F6").Value2
$g6 = $ws.Range("G6").Value2
$h6 = $ws.Range("H6").Value2
$s6 = $ws.Range("S6").Value2

$ac6 = $ws.Range("AC6").Value2
$ai6 = $ws.Range("AI6").Value2

# --- Write row 6's old values into row 4 ---
$ws.Range("A4").Value2 = $a6
$ws.Range("B4").Value2 = $b6
$ws.Range("D4").Value2 = $d6
$ws.Range("E4").Value2 = $e6
$ws.Range("F4").Value2 = $f6
$ws.Range("G4").Value2 = $g6
$ws.Range("H4").Value2 = $h6
$ws.Range("S4").Value2 = $s6

$ws.Range("AC4").Value2 = $ac6
$ws.Range("AI4").Value2 = $ai6

# --- Write row 4's old values into row 6 ---
$ws.Range("A6").Value2 = $a4
$ws.Range("B6").Value2 = $b4
$ws.Range("D6").Value2 = $d4
$ws.Range("E6").Value2 = $e4
$ws.Range("F6").Value2 = $f4
$ws.Range("G6").Value2 = $g4
$ws.Range("H6").Value2 = $h4
$ws.Range("S6").Value2 = $s4

# --- Clear the comments that moved away from row 6 ---
$ws.Range("AC6").ClearContents()
$ws.Range("AI6").ClearContents()
